$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "327.34"
Set-TextCell $ws "E2" "-0.79%"
Set-TextCell $ws "D3" "43.97"
Set-TextCell $ws "E3" "0.89%"
Set-TextCell $ws "D4" "5.562"
Set-TextCell $ws "E4" "-0.87%"
Set-TextCell $ws "D5" "0.08039"
Set-TextCell $ws "E5" "-1.95%"
Set-TextCell $ws "D6" "1.901"
Set-TextCell $ws "E6" "-0.52%"
Set-TextCell $ws "D7" "4.279"
Set-TextCell $ws "E7" "-2.83%"
Set-TextCell $ws "D8" "2.541"
Set-TextCell $ws "E8" "-10.72%"
Set-TextCell $ws "D9" "0.9444"
Set-TextCell $ws "E9" "0.06%"
Set-TextCell $ws "D10" "0.1164"
Set-TextCell $ws "E10" "-2.97%"
Set-TextCell $ws "D11" "0.1841"
Set-TextCell $ws "E11" "-4.04%"
Set-TextCell $ws "D12" "0.09688"
Set-TextCell $ws "E12" "-1.25%"
Set-TextCell $ws "D13" "0.04375"
Set-TextCell $ws "E13" "1.11%"
Set-TextCell $ws "D14" "0.1068"
Set-TextCell $ws "E14" "0.01%"
Set-TextCell $ws "D15" "0.001275"
Set-TextCell $ws "E15" "-1.50%"
Set-TextCell $ws "D16" "0.005962"
Set-TextCell $ws "E16" "-0.95%"
Set-TextCell $ws "E17" "-2.90%"
Set-TextCell $ws "D18" "0.3451"
Set-TextCell $ws "E18" "-2.42%"
Set-TextCell $ws "D19" "9.981"
Set-TextCell $ws "E19" "14.24%"
Set-TextCell $ws "D20" "0.1380"
Set-TextCell $ws "E20" "0.83%"
Set-TextCell $ws "D21" "0.2510"
Set-TextCell $ws "D22" "0.04201"
Set-TextCell $ws "E22" "-4.27%"
Set-TextCell $ws "D23" "0.001247"
Set-TextCell $ws "E23" "0.68%"
Set-TextCell $ws "D24" "0.004280"
Set-TextCell $ws "E24" "-0.69%"
Set-TextCell $ws "D25" "0.0001262"
Set-TextCell $ws "E25" "2.11%"
Set-TextCell $ws "D26" "0.0003998"
Set-TextCell $ws "E26" "-0.12%"
Set-TextCell $ws "D38" "0.02642"
Set-TextCell $ws "E38" "-5.46%"
Set-TextCell $ws "D39" "0.05500"
Set-TextCell $ws "E39" "-4.05%"
Set-TextCell $ws "D40" "0.007597"
Set-TextCell $ws "E40" "-4.30%"
Set-TextCell $ws "D41" "0.1395"
Set-TextCell $ws "E41" "-1.56%"
Set-TextCell $ws "D42" "0.008027"
Set-TextCell $ws "E42" "-17.64%"
Set-TextCell $ws "D43" "0.002012"
Set-TextCell $ws "E43" "-4.39%"
Set-TextCell $ws "D44" "0.008844"
Set-TextCell $ws "E44" "-8.52%"
Set-TextCell $ws "D45" "0.00006918"
Set-TextCell $ws "E45" "-5.88%"
Set-TextCell $ws "D46" "0.00000000752"
Set-TextCell $ws "E46" "-0.11%"
Set-TextCell $ws "D47" "0.002275"
Set-TextCell $ws "E47" "-0.12%"
Set-TextCell $ws "D48" "0.005776"
Set-TextCell $ws "E48" "66.99%"
Set-TextCell $ws "D49" "0.00002105"
Set-TextCell $ws "E49" "-0.11%"
Set-TextCell $ws "D50" "0.0002005"
Set-TextCell $ws "E50" "-0.11%"
